$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 729.63635
$ws.Cells.Item(2, 9).Value = 671.6667
$ws.Cells.Item(2, 11).Value = 671.6667
$ws.Cells.Item(2, 13).Value = -558.6667
$ws.Cells.Item(5, 8).Value = 548.8461
$ws.Cells.Item(5, 9).Value = 544.8182
$ws.Cells.Item(5, 10).Value = 571
$ws.Cells.Item(5, 11).Value = 544.8182
$ws.Cells.Item(5, 12).Value = 571
$ws.Cells.Item(5, 13).Value = -429.8182
$ws.Cells.Item(5, 14).Value = -801
$ws.Cells.Item(64, 8).Value = 11574.5
$ws.Cells.Item(64, 9).Value = 11499
$ws.Cells.Item(64, 10).Value = 11650
$ws.Cells.Item(64, 11).Value = 11499
$ws.Cells.Item(64, 12).Value = 11650
$ws.Cells.Item(64, 13).Value = -11251
$ws.Cells.Item(64, 14).Value = -12146
$ws.Cells.Item(67, 8).Value = 11574.5
$ws.Cells.Item(67, 9).Value = 11499
$ws.Cells.Item(67, 10).Value = 11650
$ws.Cells.Item(67, 11).Value = 11499
$ws.Cells.Item(67, 12).Value = 11650
$ws.Cells.Item(67, 13).Value = -10641
$ws.Cells.Item(67, 14).Value = -13366
$ws.Cells.Item(106, 8).Value = 3701.625
$ws.Cells.Item(106, 9).Value = 3572.077
$ws.Cells.Item(106, 11).Value = 3572.077
$ws.Cells.Item(106, 13).Value = -2941.077
$ws.Cells.Item(118, 8).Value = 1498.6364
$ws.Cells.Item(118, 9).Value = 1499.2222
$ws.Cells.Item(118, 11).Value = 4497.6666
$ws.Cells.Item(118, 13).Value = -2840.6666
$ws.Cells.Item(137, 8).Value = 2130.9666
$ws.Cells.Item(137, 9).Value = 2083.9565
$ws.Cells.Item(137, 11).Value = 6251.869499999999
$ws.Cells.Item(137, 13).Value = -3701.869499999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 171
$ws.Cells.Item(4, 9).Value = 160.66667
$ws.Cells.Item(4, 11).Value = 160.66667
$ws.Cells.Item(4, 13).Value = -44.66667000000001
$ws.Cells.Item(6, 8).Value = 668500.7
$ws.Cells.Item(6, 9).Value = 2751
$ws.Cells.Item(6, 11).Value = 2751
$ws.Cells.Item(6, 13).Value = -2578
$ws.Cells.Item(32, 8).Value = 15394404
$ws.Cells.Item(32, 9).Value = 23814004
$ws.Cells.Item(32, 10).Value = 19482.479
$ws.Cells.Item(32, 11).Value = 23814004
$ws.Cells.Item(32, 12).Value = 19482.479
$ws.Cells.Item(32, 13).Value = -23813717
$ws.Cells.Item(32, 14).Value = -20056.479
$ws.Cells.Item(45, 8).Value = 1347
$ws.Cells.Item(45, 9).Value = 1038.2
$ws.Cells.Item(45, 10).Value = 2119
$ws.Cells.Item(45, 11).Value = 1038.2
$ws.Cells.Item(45, 12).Value = 2119
$ws.Cells.Item(45, 13).Value = -661.2
$ws.Cells.Item(45, 14).Value = -2873
$ws.Cells.Item(61, 8).Value = 25644642
$ws.Cells.Item(61, 9).Value = 34484890
$ws.Cells.Item(61, 10).Value = 7919.8
$ws.Cells.Item(61, 11).Value = 34484890
$ws.Cells.Item(61, 12).Value = 7919.8
$ws.Cells.Item(61, 13).Value = -34484678
$ws.Cells.Item(61, 14).Value = -8343.799999999999
$ws.Cells.Item(110, 8).Value = 14921
$ws.Cells.Item(110, 9).Value = 16813.738
$ws.Cells.Item(110, 11).Value = 16813.738
$ws.Cells.Item(110, 13).Value = -14768.738
$ws.Cells.Item(132, 8).Value = 47631976
$ws.Cells.Item(132, 9).Value = 14821.647
$ws.Cells.Item(132, 11).Value = 44464.94100000001
$ws.Cells.Item(132, 13).Value = -41934.94100000001
$ws.Cells.Item(136, 8).Value = 25644642
$ws.Cells.Item(136, 9).Value = 34484890
$ws.Cells.Item(136, 10).Value = 7919.8
$ws.Cells.Item(136, 11).Value = 103454670
$ws.Cells.Item(136, 12).Value = 23759.4
$ws.Cells.Item(136, 13).Value = -103452120
$ws.Cells.Item(136, 14).Value = -28859.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2733.15
$ws.Cells.Item(134, 9).Value = 2760.7715
$ws.Cells.Item(134, 11).Value = 8282.3145
$ws.Cells.Item(134, 13).Value = -5747.3145
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 392.57144
$ws.Cells.Item(7, 9).Value = 167
$ws.Cells.Item(7, 11).Value = 167
$ws.Cells.Item(7, 13).Value = -54
$ws.Cells.Item(31, 8).Value = 21283432
$ws.Cells.Item(31, 9).Value = 6410.0625
$ws.Cells.Item(31, 10).Value = 66674410
$ws.Cells.Item(31, 11).Value = 6410.0625
$ws.Cells.Item(31, 12).Value = 66674410
$ws.Cells.Item(31, 13).Value = -6115.0625
$ws.Cells.Item(31, 14).Value = -66675000
$ws.Cells.Item(34, 8).Value = 21283432
$ws.Cells.Item(34, 9).Value = 6410.0625
$ws.Cells.Item(34, 10).Value = 66674410
$ws.Cells.Item(34, 11).Value = 6410.0625
$ws.Cells.Item(34, 12).Value = 66674410
$ws.Cells.Item(34, 13).Value = -6208.0625
$ws.Cells.Item(34, 14).Value = -66674814
$ws.Cells.Item(141, 8).Value = 202862.14
$ws.Cells.Item(141, 10).Value = 281575.78
$ws.Cells.Item(141, 12).Value = 281575.78
$ws.Cells.Item(141, 14).Value = -291935.78
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 3022.1
$ws.Cells.Item(75, 9).Value = 2518.8333
$ws.Cells.Item(75, 10).Value = 3777
$ws.Cells.Item(75, 11).Value = 7556.499899999999
$ws.Cells.Item(75, 12).Value = 11331
$ws.Cells.Item(75, 13).Value = -6558.499899999999
$ws.Cells.Item(75, 14).Value = -13327
$ws.Cells.Item(78, 8).Value = 3022.1
$ws.Cells.Item(78, 9).Value = 2518.8333
$ws.Cells.Item(78, 10).Value = 3777
$ws.Cells.Item(78, 11).Value = 22669.4997
$ws.Cells.Item(78, 12).Value = 33993
$ws.Cells.Item(78, 13).Value = -17677.4997
$ws.Cells.Item(78, 14).Value = -43977
$ws.Cells.Item(80, 8).Value = 3251.5
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(83, 8).Value = 3251.5
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(113, 8).Value = 2505.318
$ws.Cells.Item(113, 9).Value = 1550.2222
$ws.Cells.Item(113, 10).Value = 3166.5386
$ws.Cells.Item(113, 11).Value = 4650.6666
$ws.Cells.Item(113, 12).Value = 9499.6158
$ws.Cells.Item(113, 13).Value = -2480.6666
$ws.Cells.Item(113, 14).Value = -13839.6158
$ws.Cells.Item(132, 8).Value = 3336921.8
$ws.Cells.Item(132, 9).Value = 2559.4
$ws.Cells.Item(132, 10).Value = 4448376
$ws.Cells.Item(132, 11).Value = 23034.6
$ws.Cells.Item(132, 12).Value = 40035384
$ws.Cells.Item(132, 13).Value = -20504.6
$ws.Cells.Item(132, 14).Value = -40040444
$ws.Cells.Item(134, 8).Value = 8017.579
$ws.Cells.Item(134, 9).Value = 2718
$ws.Cells.Item(134, 10).Value = 19500
$ws.Cells.Item(134, 11).Value = 8154
$ws.Cells.Item(134, 12).Value = 58500
$ws.Cells.Item(134, 13).Value = -3084
$ws.Cells.Item(134, 14).Value = -68640
$ws.Cells.Item(136, 8).Value = 5694
$ws.Cells.Item(136, 10).Value = 7225
$ws.Cells.Item(136, 12).Value = 21675
$ws.Cells.Item(136, 14).Value = -31875
$ws.Cells.Item(137, 8).Value = 5235.4
$ws.Cells.Item(137, 9).Value = 2794.25
$ws.Cells.Item(137, 10).Value = 15000
$ws.Cells.Item(137, 11).Value = 8382.75
$ws.Cells.Item(137, 12).Value = 45000
$ws.Cells.Item(137, 13).Value = -3282.75
$ws.Cells.Item(137, 14).Value = -55200
$ws.Cells.Item(138, 8).Value = 1757.6666
$ws.Cells.Item(138, 9).Value = 1120
$ws.Cells.Item(138, 11).Value = 3360
$ws.Cells.Item(138, 13).Value = 1780
$ws.Cells.Item(139, 8).Value = 2679.5
$ws.Cells.Item(139, 9).Value = 2608.8
$ws.Cells.Item(139, 10).Value = 3033
$ws.Cells.Item(139, 11).Value = 7826.400000000001
$ws.Cells.Item(139, 12).Value = 9099
$ws.Cells.Item(139, 13).Value = -2686.400000000001
$ws.Cells.Item(139, 14).Value = -19379
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 17299
$ws.Cells.Item(43, 9).Value = 2900
$ws.Cells.Item(43, 10).Value = 24498.5
$ws.Cells.Item(43, 11).Value = 2900
$ws.Cells.Item(43, 12).Value = 24498.5
$ws.Cells.Item(43, 13).Value = -2749
$ws.Cells.Item(43, 14).Value = -24800.5
$ws.Cells.Item(113, 8).Value = 2275.2058
$ws.Cells.Item(113, 10).Value = 3181.4666
$ws.Cells.Item(113, 12).Value = 3181.4666
$ws.Cells.Item(113, 14).Value = -7521.4666
$ws.Cells.Item(132, 8).Value = 2212.625
$ws.Cells.Item(132, 9).Value = 2481.72
$ws.Cells.Item(132, 11).Value = 7445.16
$ws.Cells.Item(132, 13).Value = -4915.16
$ws.Cells.Item(139, 8).Value = 123030.14
$ws.Cells.Item(139, 10).Value = 123030.14
$ws.Cells.Item(139, 12).Value = 123030.14
$ws.Cells.Item(139, 14).Value = -133310.14
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 5538
$ws.Cells.Item(61, 9).Value = 4999
$ws.Cells.Item(61, 11).Value = 4999
$ws.Cells.Item(61, 13).Value = -4797
$ws.Cells.Item(74, 8).Value = 62996.375
$ws.Cells.Item(74, 9).Value = 61959
$ws.Cells.Item(74, 10).Value = 66108.5
$ws.Cells.Item(74, 11).Value = 61959
$ws.Cells.Item(74, 12).Value = 66108.5
$ws.Cells.Item(74, 13).Value = -60961
$ws.Cells.Item(74, 14).Value = -68104.5
$ws.Cells.Item(77, 8).Value = 62996.375
$ws.Cells.Item(77, 9).Value = 61959
$ws.Cells.Item(77, 10).Value = 66108.5
$ws.Cells.Item(77, 11).Value = 185877
$ws.Cells.Item(77, 12).Value = 198325.5
$ws.Cells.Item(77, 13).Value = -180885
$ws.Cells.Item(77, 14).Value = -208309.5
$ws.Cells.Item(113, 8).Value = 5538
$ws.Cells.Item(113, 9).Value = 4999
$ws.Cells.Item(113, 11).Value = 4999
$ws.Cells.Item(113, 13).Value = -2829
$ws.Cells.Item(122, 8).Value = 4064.5293
$ws.Cells.Item(122, 9).Value = 3632
$ws.Cells.Item(122, 11).Value = 10896
$ws.Cells.Item(122, 13).Value = -8446
$ws.Cells.Item(131, 8).Value = 89078
$ws.Cells.Item(131, 10).Value = 89078
$ws.Cells.Item(131, 12).Value = 89078
$ws.Cells.Item(131, 14).Value = -99158
$ws.Cells.Item(132, 8).Value = 62501304
$ws.Cells.Item(132, 9).Value = 1288.4445
$ws.Cells.Item(132, 10).Value = 400001400
$ws.Cells.Item(132, 11).Value = 3865.3335
$ws.Cells.Item(132, 12).Value = 1200004200
$ws.Cells.Item(132, 13).Value = -1335.3335
$ws.Cells.Item(132, 14).Value = -1200009260
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 38137
$ws.Cells.Item(92, 10).Value = 38137
$ws.Cells.Item(92, 12).Value = 38137
$ws.Cells.Item(92, 14).Value = -43129
$ws.Cells.Item(136, 8).Value = 1346.381
$ws.Cells.Item(136, 9).Value = 1200.9736
$ws.Cells.Item(136, 11).Value = 3602.9208
$ws.Cells.Item(136, 13).Value = -1052.9208
